$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-02-17 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-02-18 Sunday", 2)

$d.Content.Find.Execute("508÷3=169, 1", $true, $false, $false, $false, $false, $true, 1, $false, "889÷2=444, 1", 2)
$d.Content.Find.Execute("789÷6=131, 3", $true, $false, $false, $false, $false, $true, 1, $false, "901÷3=300, 1", 2)
$d.Content.Find.Execute("762÷3=254, 0", $true, $false, $false, $false, $false, $true, 1, $false, "360÷3=120, 0", 2)
$d.Content.Find.Execute("714÷2=357, 0", $true, $false, $false, $false, $false, $true, 1, $false, "229÷4=57, 1", 2)
$d.Content.Find.Execute("612÷2=306, 0", $true, $false, $false, $false, $false, $true, 1, $false, "471÷9=52, 3", 2)

$d.Content.Find.Execute("594÷3=198, 0", $true, $false, $false, $false, $false, $true, 1, $false, "590÷6=98, 2", 2)
$d.Content.Find.Execute("526÷2=263, 0", $true, $false, $false, $false, $false, $true, 1, $false, "148÷8=18, 4", 2)
$d.Content.Find.Execute("699÷5=139, 4", $true, $false, $false, $false, $false, $true, 1, $false, "821÷3=273, 2", 2)
$d.Content.Find.Execute("818÷6=136, 2", $true, $false, $false, $false, $false, $true, 1, $false, "608÷9=67, 5", 2)
$d.Content.Find.Execute("641÷9=71, 2", $true, $false, $false, $false, $false, $true, 1, $false, "876÷5=175, 1", 2)

$d.Content.Find.Execute("317÷2=158, 1", $true, $false, $false, $false, $false, $true, 1, $false, "101÷8=12, 5", 2)
$d.Content.Find.Execute("169÷7=24, 1", $true, $false, $false, $false, $false, $true, 1, $false, "198÷7=28, 2", 2)
$d.Content.Find.Execute("816÷9=90, 6", $true, $false, $false, $false, $false, $true, 1, $false, "108÷8=13, 4", 2)
$d.Content.Find.Execute("605÷6=100, 5", $true, $false, $false, $false, $false, $true, 1, $false, "541÷5=108, 1", 2)
$d.Content.Find.Execute("665÷9=73, 8", $true, $false, $false, $false, $false, $true, 1, $false, "213÷7=30, 3", 2)

$d.Content.Find.Execute("408÷2=204, 0", $true, $false, $false, $false, $false, $true, 1, $false, "823÷6=137, 1", 2)
$d.Content.Find.Execute("915÷4=228, 3", $true, $false, $false, $false, $false, $true, 1, $false, "175÷9=19, 4", 2)
$d.Content.Find.Execute("146÷3=48, 2", $true, $false, $false, $false, $false, $true, 1, $false, "271÷5=54, 1", 2)
$d.Content.Find.Execute("563÷4=140, 3", $true, $false, $false, $false, $false, $true, 1, $false, "860÷4=215, 0", 2)
$d.Content.Find.Execute("133÷2=66, 1", $true, $false, $false, $false, $false, $true, 1, $false, "452÷6=75, 2", 2)

$d.Content.Find.Execute("659÷8=82, 3", $true, $false, $false, $false, $false, $true, 1, $false, "493÷6=82, 1", 2)
$d.Content.Find.Execute("751÷3=250, 1", $true, $false, $false, $false, $false, $true, 1, $false, "698÷3=232, 2", 2)
$d.Content.Find.Execute("738÷6=123, 0", $true, $false, $false, $false, $false, $true, 1, $false, "173÷8=21, 5", 2)
$d.Content.Find.Execute("791÷5=158, 1", $true, $false, $false, $false, $false, $true, 1, $false, "321÷7=45, 6", 2)
$d.Content.Find.Execute("146÷5=29, 1", $true, $false, $false, $false, $false, $true, 1, $false, "600÷8=75, 0", 2)
